# Update the "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" sheets to reflect a newer data scrape.
# Note: the F8 value differs slightly between the two sheets in the
# source data (7752 vs 7753), so they are set independently.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 8
$ws1.Range("F7").Value  = 294
$ws1.Range("F8").Value  = 7752
$ws1.Range("F9").Value  = 72
$ws1.Range("F13").Value = 4
$ws1.Range("F15").Value = 18
$ws1.Range("F19").Value = 681
$ws1.Range("F20").Value = 19

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 8
$ws4.Range("F7").Value  = 294
$ws4.Range("F8").Value  = 7753
$ws4.Range("F9").Value  = 72
$ws4.Range("F13").Value = 4
$ws4.Range("F15").Value = 18
$ws4.Range("F19").Value = 681
$ws4.Range("F20").Value = 19
